$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 81

$ws.Cells.Item($newRow, 1).Value = "2025-04-29 12:00:52"
$ws.Cells.Item($newRow, 2).Value = 228
